# Refresh the cryptos price/volume table with the latest scraped values.
# (GitHub Actions crypto-price refresh run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.332.09'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.810.64'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = "'313.64"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').Value = "'0.9998"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').Value = "'0.5143"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('D8').Value = "'0.3994"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.18%  '
$ws.Range('D9').Value = "'0.07867"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.39%  '
$ws.Range('D10').Value = "'1.114"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').Value = "'40.86"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.56%  '
$ws.Range('D12').Value = "'6.384"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = "'1.000"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = "'20.39"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.05%  '
$ws.Range('D15').Value = "'7.362"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').Value = '1.802.81'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').Value = "'92.96"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('E18').Value = '  -3.67%  '
$ws.Range('D19').Value = "'0.06575"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.18%  '
$ws.Range('D20').Value = "'0.9994"
$ws.Range('D20').ClearFormats()
$ws.Range('D21').Value = "'17.35"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('D22').Value = "'6.029"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('D23').Value = '28.385.53'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').Value = "'2.225"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('D26').Value = "'160.78"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('E27').Value = '  -2.99%  '
$ws.Range('D28').Value = '2.016.24'
$ws.Range('D29').Value = "'2.402"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').Value = "'128.78"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D31').Value = "'0.1087"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('D32').Value = "'1.065"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.39%  '
$ws.Range('D33').Value = "'3.664"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').Value = "'5.585"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.75%  '
$ws.Range('D35').Value = "'0.07239"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('D36').Value = "'9.185"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.33%  '
$ws.Range('D37').Value = "'0.02342"
$ws.Range('D37').ClearFormats()
$ws.Range('E38').Value = '  -2.60%  '
$ws.Range('D39').Value = "'5.080"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.39%  '
$ws.Range('D40').Value = "'11.62"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('D41').Value = "'0.6205"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.17%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.162"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.89%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').Value = "'0.9991"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = "'0.6013"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.46%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'13.18"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.93%  '
$ws.Range('D46').Value = "'1.311"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.09%  '
$ws.Range('D47').Value = "'3.743"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('D48').Value = "'125.80"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('D49').Value = "'1.224"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').Value = "'1.937"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.27%  '
$ws.Range('D51').Value = "'0.06850"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.96%  '
